# Updated cryptos list - applies changed cell values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.017.44'
$ws.Range("E2").Value = '  -2.68%  '
$ws.Range("D3").Value = '2.639.37'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'527.34"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = "'144.49"
$ws.Range("E6").Value = '  -3.82%  '
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = "'0.571"
$ws.Range("E8").Value = '  -1.30%  '
$ws.Range("D9").Value = "'6.64"
$ws.Range("E9").Value = '  -5.74%  '
$ws.Range("D10").Value = "'0.104"
$ws.Range("E10").Value = '  -1.54%  '
$ws.Range("D11").Value = "'0.337"
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").Value = '3.101.81'
$ws.Range("E13").Value = '  -2.38%  '
$ws.Range("D14").Value = '58.990.68'
$ws.Range("E14").Value = '  -2.69%  '
$ws.Range("D15").Value = "'21.00"
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("E16").Value = '  -1.23%  '
$ws.Range("D17").Value = '2.594.57'
$ws.Range("E17").Value = '  -4.24%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").Value = "'341.52"
$ws.Range("E18").Value = '  -1.33%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = "'4.46"
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("D20").Value = "'10.57"
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Value = "'6.34"
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = "'65.51"
$ws.Range("E23").Value = '  +2.94%  '
$ws.Range("D24").Value = "'0.419"
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("E25").Value = '  -1.95%  '
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("D27").Value = "'7.25"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").Value = '  -3.23%  '
$ws.Range("E29").Value = '  -4.59%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("D32").Value = "'18.94"
$ws.Range("E32").Value = '  -0.75%  '
$ws.Range("D33").Value = "'149.96"
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("D34").Value = "'4.24"
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("B35").Value = 'SuiNetwork'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D35").Value = "'0.932"
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = "'1.20"
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("D37").Value = "'0.875"
$ws.Range("E37").Value = '  -2.98%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = "'1.48"
$ws.Range("E38").Value = '  -3.47%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = "'36.57"
$ws.Range("E39").Value = '  -2.13%  '
$ws.Range("D40").Value = "'3.66"
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = "'0.605"
$ws.Range("E42").Value = '  -5.45%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = "'0.0976"
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("D44").Value = "'271.23"
$ws.Range("E44").Value = '  -2.85%  '
$ws.Range("D45").Value = "'19.44"
$ws.Range("E45").Value = '  -3.79%  '
$ws.Range("D46").Value = "'0.0539"
$ws.Range("E46").Value = '  -1.05%  '
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("D48").Value = '2.050.37'
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = "'4.80"
$ws.Range("E49").Value = '  -3.82%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = "'0.0230"
$ws.Range("E50").Value = '  -1.47%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = "'19.04"
$ws.Range("E51").Value = '  -2.29%  '
